$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.018.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.779.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.87"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.780.08"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.54"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.73"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.413.74"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.781.27"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.83"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.943.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.82"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -8.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.76"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.16"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.928.23"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.61"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.59"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.26"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.736.95"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.73"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "402.73"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.64"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -7.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.43"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.96"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.28%  "
